$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.562.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.35%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.911.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.48%  '

$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '

$ws.Range('E6').Value = '  -0.31%  '

$ws.Range('E7').Value = '  +1.72%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3961'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09753'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.06%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.159'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.41%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.37%  '

$ws.Range('E12').Value = '  +2.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.907.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.63%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.579'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.19%  '

$ws.Range('E16').Value = '  -0.23%  '

$ws.Range('E17').Value = '  +1.83%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.75%  '

$ws.Range('E19').Value = '  -0.05%  '

$ws.Range('E20').Value = '  +5.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.296'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.53%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.623.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.47'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.281'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.745'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +15.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.127.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.40%  '

$ws.Range('E28').Value = '  +3.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '159.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.104'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1073'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.66%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.741'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.646'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.72%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.898'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.90%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06811'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.49%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02447'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.29%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.273'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.26%  '

$ws.Range('E39').Value = '  +4.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.74'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.88%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.113'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6440'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.191'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.12%  '

$ws.Range('E44').Value = '  -0.31%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.66'
$ws.Range('D45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6108'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.98%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.286'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.51%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.671'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.77%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.046'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.09%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.85%  '

$ws.Range('E51').Value = '  +3.20%  '
